# Update gh-pages output (南宁-漫展信息.xlsx)
# - bumps "想去人数" (F column) counters on the 展览 (sheet 1) and 全部类型 (sheet 4) tabs
# - appends a new row for 南宁·万圣漫控嘉年华10 at the bottom of each of those tabs

$wb = $excel.ActiveWorkbook

$sheetExhibit = $wb.Worksheets.Item(1)   # 展览
$sheetAll     = $wb.Worksheets.Item(4)   # 全部类型

# ---- 展览 (sheet 1): refresh "想去人数" counts ----
$sheetExhibit.Range("F2").Value  = 69
$sheetExhibit.Range("F3").Value  = 548
$sheetExhibit.Range("F5").Value  = 288
$sheetExhibit.Range("F6").Value  = 394
$sheetExhibit.Range("F8").Value  = 2319
$sheetExhibit.Range("F10").Value = 5780
$sheetExhibit.Range("F12").Value = 378

# ---- 全部类型 (sheet 4): refresh "想去人数" counts ----
$sheetAll.Range("F2").Value  = 69
$sheetAll.Range("F4").Value  = 548
$sheetAll.Range("F6").Value  = 288
$sheetAll.Range("F7").Value  = 394
$sheetAll.Range("F11").Value = 2319
$sheetAll.Range("F13").Value = 5780
$sheetAll.Range("F15").Value = 378

# ---- append new row: 南宁·万圣漫控嘉年华10 ----

function Add-ConEntry($ws, $rowNum, $idValue) {
    $prevRow = $rowNum - 1

    $ws.Range("A$rowNum").Value = $idValue

    # Force the start-date cell to stay text (like the rest of column B)
    # instead of being auto-parsed into a date serial, then restore the
    # plain/default formatting used by the sibling date cells above it.
    $ws.Range("B$rowNum").NumberFormat = "@"
    $ws.Range("B$rowNum").Value = "2024-11-02"
    $ws.Range("B$prevRow").Copy()
    $ws.Range("B$rowNum").PasteSpecial(-4122)

    $ws.Range("C$rowNum").Value = "南宁·万圣漫控嘉年华10"
    $ws.Range("D$rowNum").Value = "亭洪路45号 百益上河城"
    $ws.Range("E$rowNum").Value = "2024.11.02 11:00-11.03 22:00"
    $ws.Range("F$rowNum").Value = 1
    $ws.Range("G$rowNum").Value = 50
    $ws.Range("H$rowNum").Value = "https://show.bilibili.com/platform/detail.html?id=87820"
    $ws.Range("I$rowNum").Value = "//i1.hdslb.com/bfs/openplatform/202406/abJD2cvV1718955681653.jpeg"

    # Match the bold/centered/bordered style used by the rest of column A
    $ws.Range("A$prevRow").Copy()
    $ws.Range("A$rowNum").PasteSpecial(-4122)
}

Add-ConEntry $sheetExhibit 13 12
Add-ConEntry $sheetAll 17 16

$excel.CutCopyMode = 0
